# Generate Report for Handback
#
# The handback for the de-de / zh-cn localization files is now in sync
# with en-US, so refresh the generated status report:
#   - flip the status text from "Ready for handoff" to the handed-back state
#   - stamp the latest handback datetime for each locale
#   - clear the stale "version not latest" error now that things are in sync
#   - let column widths track the new (longer/shorter) cell contents

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-08-29 06:50:29"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-29 06:50:36"
$wsDeDe.Range("P2").Value = ""

# Resize the columns that now hold different-length text, same as Excel
# does when the report is regenerated.
$wsOverview.Columns.Item(5).ColumnWidth = 29.16
$wsOverview.Columns.Item(6).ColumnWidth = 29.16

$wsZhCn.Columns.Item(3).ColumnWidth = 29.16
$wsZhCn.Columns.Item(16).ColumnWidth = 12.83

$wsDeDe.Columns.Item(3).ColumnWidth = 29.16
$wsDeDe.Columns.Item(16).ColumnWidth = 12.83
